# Append a new "2021年" data row (row 5) to Sheet1, mirroring the layout of
# the existing year rows (2018-2020 in rows 2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous year row (row 4) down onto the new
# row 5 first, so the year-label cell (column A) picks up the same style
# (bold, bordered, centered) as A2:A4 without touching B:AT (which carry no
# explicit style).
$ws.Range("A4:AT4").Copy()
$ws.Range("A5:AT5").PasteSpecial(-4122)  # xlPasteFormats

# Year label
$ws.Range("A5").Value = "2021年"

# Data values for 2021, in the same column order as the header row (row 1)
$ws.Range("B5").Value = 53.038
$ws.Range("C5").Value = 58.282
$ws.Range("D5").Value = 38.544
$ws.Range("E5").Value = 41.624
$ws.Range("F5").Value = 20
$ws.Range("G5").Value = 27.536
$ws.Range("H5").Value = 40.12
$ws.Range("I5").Value = 44.48
$ws.Range("J5").Value = 46.113
$ws.Range("K5").Value = 53.841
$ws.Range("L5").Value = 38.983
$ws.Range("M5").Value = 34.256
$ws.Range("N5").Value = 29.024
$ws.Range("O5").Value = 18.553
$ws.Range("P5").Value = 34.748
$ws.Range("Q5").Value = 38.146
$ws.Range("R5").Value = 22.165
$ws.Range("S5").Value = 20.282
$ws.Range("T5").Value = 40.855
$ws.Range("U5").Value = 17.474
$ws.Range("V5").Value = 49.36
$ws.Range("W5").Value = 59.524
$ws.Range("X5").Value = 15.688
$ws.Range("Y5").Value = 13.452
$ws.Range("Z5").Value = 15.151
$ws.Range("AA5").Value = 14.954
$ws.Range("AB5").Value = 50.449
$ws.Range("AC5").Value = 25.838
$ws.Range("AD5").Value = 33.766
$ws.Range("AE5").Value = 37.762
$ws.Range("AF5").Value = 34.713
$ws.Range("AG5").Value = 21.754
$ws.Range("AH5").Value = 56.024
$ws.Range("AI5").Value = 50.156
$ws.Range("AJ5").Value = 34.826
$ws.Range("AK5").Value = 33.286
$ws.Range("AL5").Value = 16.927
$ws.Range("AM5").Value = 34.639
$ws.Range("AN5").Value = 38.603
$ws.Range("AO5").Value = 47.015
$ws.Range("AP5").Value = 27.828
$ws.Range("AQ5").Value = 15.943
$ws.Range("AR5").Value = 38.313
$ws.Range("AS5").Value = 32.371
$ws.Range("AT5").Value = 16.689
